# "revised based on class notes"
# Slide 2, Content Placeholder, paragraph "Exceptions & asserts":
# merge the separate "&" and " asserts" runs into a single "& asserts" run.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame2.TextRange

$para = $tr.Paragraphs(4, 1)

# Run 3 is "&", run 4 is " asserts" -- fold run 4's text into run 3,
# then clear run 4 so it collapses out of the run list.
$ampRun = $para.Runs(3, 1)
$ampRun.Text = "& asserts"

$tailRun = $para.Runs(4, 1)
$tailRun.Text = ""
